$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 678.21
$ws.Range("C3").Value = 685.6799999999999
$ws.Range("C4").Value = 670.89
$ws.Range("C5").Value = 667.78
$ws.Range("C6").Value = 667.78
